$wb = $excel.ActiveWorkbook

# --- Sheet "Means": update column F (Within 5 miles of HFC production facility), rows 2-10 ---
$wsMeans = $wb.Worksheets.Item("Means")
$wsMeans.Range("F2").Value = 82
$wsMeans.Range("F3").Value = 15
$wsMeans.Range("F4").Value = 2.9
$wsMeans.Range("F5").Value = 3.4
$wsMeans.Range("F6").Value = 54
$wsMeans.Range("F7").Value = 11
$wsMeans.Range("F8").Value = 4.2
$wsMeans.Range("F9").Value = 50
$wsMeans.Range("F10").Value = 0.57

# --- Sheet "Standard Deviations": update column F, rows 2-10 (except row 7, unchanged) ---
$wsSD = $wb.Worksheets.Item("Standard Deviations")
$wsSD.Range("F2").Value = 15
$wsSD.Range("F3").Value = 17
$wsSD.Range("F4").Value = 4.8
$wsSD.Range("F5").Value = 2.2
$wsSD.Range("F6").Value = 30
$wsSD.Range("F8").Value = 2.3
$wsSD.Range("F9").Value = 0
$wsSD.Range("F10").Value = 0.052
